$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$v = $ws.Range("C2").Value
Write-Host "Cell C2 value:" $v
$v2 = $ws.Range("C2").Value2
Write-Host "Cell C2 value2:" $v2
